$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove C2 entirely, update E2
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 9.591339540850829

# Row 4
$ws.Range("E4").Value = -2.225127715916664

# Row 5
$ws.Range("C5").Value = 9.399485634179205

# Row 6
$ws.Range("C6").Value = 5.169490031659651
$ws.Range("E6").Value = 9.213376886330327

# Row 7
$ws.Range("C7").Value = -0.3722371047999995
$ws.Range("E7").Value = 2.684220738731979

# Row 8
$ws.Range("C8").Value = 4.098801479368319

# Row 9
$ws.Range("E9").Value = 2.714258593289998

# Row 10
$ws.Range("C10").Value = 2.352205130086094

# Row 11
$ws.Range("C11").Value = 4.083548352538391
$ws.Range("E11").Value = 3.58625614607444

# Row 12
$ws.Range("C12").Value = 4.861590900330715

# Row 14
$ws.Range("E14").Value = -0.6322362079330235

# Row 15
$ws.Range("E15").Value = 1.077755602068309

# Row 18
$ws.Range("E18").Value = 1.793234865396331

# Row 19
$ws.Range("E19").Value = -0.6714033493142035
